$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.932809352874756
$ws.Range("B1").Value = 3.159520387649536
$ws.Range("C1").Value = 2.9838547706604
$ws.Range("D1").Value = 1.024213433265686
$ws.Range("E1").Value = 0.6654894351959229
